$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: two new observation rows are added at the top of the
# "Terminal Hortofrutícola Agro Chillán" block (rows 415-416), pushing the
# previously-existing rows 415-421 down to 417-423 unchanged.
$ws.Rows("415:416").Insert()

# New row 415
$ws.Range("A415").Value = 7
$ws.Range("B415").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C415").Value = "Ñuble"
$ws.Range("D415").Value = 44628
$ws.Range("E415").Value = 16
$ws.Range("F415").Value = 100112020
$ws.Range("G415").Value = "Tomate"
$ws.Range("H415").Value = "Larga vida"
$ws.Range("I415").Value = "Primera"
$ws.Range("J415").Value = 300
$ws.Range("K415").Value = 10000
$ws.Range("L415").Value = 10000
$ws.Range("M415").Value = 10000
$ws.Range("N415").Value = "`$/bandeja 18 kilos"
$ws.Range("O415").Value = "Región del Maule"
$ws.Range("P415").Value = 556
$ws.Range("Q415").Value = 18
$ws.Range("R415").Value = "Hortaliza"

# New row 416
$ws.Range("A416").Value = 7
$ws.Range("B416").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C416").Value = "Ñuble"
$ws.Range("D416").Value = 44628
$ws.Range("E416").Value = 16
$ws.Range("F416").Value = 100112020
$ws.Range("G416").Value = "Tomate"
$ws.Range("H416").Value = "Larga vida"
$ws.Range("I416").Value = "Primera"
$ws.Range("J416").Value = 600
$ws.Range("K416").Value = 6000
$ws.Range("L416").Value = 6500
$ws.Range("M416").Value = 6250
$ws.Range("N416").Value = "`$/caja 15 kilos"
$ws.Range("O416").Value = "Región del Maule"
$ws.Range("P416").Value = 417
$ws.Range("Q416").Value = 15
$ws.Range("R416").Value = "Hortaliza"
